$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = -0.002849999999999998
$ws.Cells.Item(2, 5).Value = -0.125
$ws.Cells.Item(2, 7).Value = 0.04860779634049324
$ws.Cells.Item(2, 8).Value = 0.04860779634049324
$ws.Cells.Item(2, 9).Value = 0.06400306618370689
$ws.Cells.Item(2, 10).Value = 0.05972805222105471
$ws.Cells.Item(2, 11).Value = 90.11
$ws.Cells.Item(2, 12).Value = 0.716865552903739
$ws.Cells.Item(2, 13).Value = 111.48
$ws.Cells.Item(2, 14).Value = 0.04229616420685207
$ws.Cells.Item(2, 15).Value = 1.237154588835867
$ws.Cells.Item(2, 16).Value = 111.48
$ws.Cells.Item(2, 17).Value = 0.04229616420685207
$ws.Cells.Item(2, 18).Value = 1.237154588835867
$ws.Cells.Item(2, 21).Value = 69.95
$ws.Cells.Item(2, 22).Value = 0.02653943923815305
$ws.Cells.Item(2, 23).Value = 0.06228932584269663
$ws.Cells.Item(2, 24).Value = 0.05927785053632793
$ws.Cells.Item(2, 25).Value = 0.003011475306368701
$ws.Cells.Item(2, 26).Value = 0.03051137447448289
$ws.Cells.Item(2, 27).Value = 0.006211126134845741
$ws.Cells.Item(2, 28).Value = 0.04751674309153629
$ws.Cells.Item(2, 29).Value = -0.04130561695669054
$ws.Cells.Item(2, 30).Value = 2669.09
$ws.Cells.Item(2, 31).Value = 7.874072903540214
$ws.Cells.Item(2, 32).Value = 2676.96407290354
$ws.Cells.Item(2, 33).Value = 2607.01407290354
$ws.Cells.Item(2, 34).Value = 0.503883557508746
$ws.Cells.Item(2, 35).Value = 0.6665142960927566
$ws.Cells.Item(2, 36).Value = 0.4972642102260812
$ws.Cells.Item(2, 37).Value = 0.6606032779995262
$ws.Cells.Item(2, 38).Value = 0.181
$ws.Cells.Item(2, 39).Value = 0.181
$ws.Cells.Item(2, 40).Value = 271.8014256619144
$ws.Cells.Item(2, 41).Value = 30.55248618784531
$ws.Cells.Item(2, 42).Value = 265.4800481571833
$ws.Cells.Item(2, 43).Value = 30.55248618784531

# Row 3
$ws.Cells.Item(3, 4).Value = -0.0294
$ws.Cells.Item(3, 5).Value = -0.125
$ws.Cells.Item(3, 7).Value = 0.2279850746268657
$ws.Cells.Item(3, 8).Value = 0.2279850746268657
$ws.Cells.Item(3, 9).Value = 0.2389876269829552
$ws.Cells.Item(3, 10).Value = 0.1910988849348439
$ws.Cells.Item(3, 11).Value = 4.19
$ws.Cells.Item(3, 12).Value = 0.1563432835820896
$ws.Cells.Item(3, 13).Value = 3.98
$ws.Cells.Item(3, 14).Value = 0.06513911620294599
$ws.Cells.Item(3, 15).Value = 0.9498806682577565
$ws.Cells.Item(3, 16).Value = 3.98
$ws.Cells.Item(3, 17).Value = 0.06513911620294599
$ws.Cells.Item(3, 18).Value = 0.9498806682577565
$ws.Cells.Item(3, 21).Value = 3.85
$ws.Cells.Item(3, 22).Value = 0.0630114566284779
$ws.Cells.Item(3, 23).Value = 0.1180281690140845
$ws.Cells.Item(3, 24).Value = 0.0464371020039641
$ws.Cells.Item(3, 25).Value = 0.07159106701012041
$ws.Cells.Item(3, 26).Value = 0.7955908123422586
$ws.Cells.Item(3, 27).Value = 0.1520365171030123
$ws.Cells.Item(3, 28).Value = 0.04467690957835233
$ws.Cells.Item(3, 29).Value = 0.10735960752466
$ws.Cells.Item(3, 30).Value = 4.19
$ws.Cells.Item(3, 31).Value = 1.775657984284003
$ws.Cells.Item(3, 32).Value = 5.965657984284004
$ws.Cells.Item(3, 33).Value = 2.115657984284004
$ws.Cells.Item(3, 34).Value = 0.08895250063276769
$ws.Cells.Item(3, 35).Value = 0.1563095803756498
$ws.Cells.Item(3, 36).Value = 0.03346730939366281
$ws.Cells.Item(3, 37).Value = 0.06165284620953325
$ws.Cells.Item(3, 38).Value = 0.181
$ws.Cells.Item(3, 39).Value = 0.181
$ws.Cells.Item(3, 40).Value = 0.6020114942528736
$ws.Cells.Item(3, 41).Value = 30.55248618784531
$ws.Cells.Item(3, 42).Value = 0.3039738483166672
$ws.Cells.Item(3, 43).Value = 30.55248618784531

# Row 4
$ws.Cells.Item(4, 4).Value = 0.0237
$ws.Cells.Item(4, 9).Value = 0.05300666353367193
$ws.Cells.Item(4, 10).Value = 0.05300666353367193
$ws.Cells.Item(4, 11).Value = -2.78
$ws.Cells.Item(4, 12).Value = -0.1139344262295082
$ws.Cells.Item(4, 21).Value = 40
$ws.Cells.Item(4, 22).Value = 0.1845869866174435
$ws.Cells.Item(4, 23).Value = -0.02195892575039494
$ws.Cells.Item(4, 24).Value = 0.05927785053632793
$ws.Cells.Item(4, 25).Value = -0.08123677628672288
$ws.Cells.Item(4, 26).Value = 0.1171763269140716
$ws.Cells.Item(4, 27).Value = 0.006211126134845741
$ws.Cells.Item(4, 28).Value = 0.04751674309153629
$ws.Cells.Item(4, 29).Value = -0.04130561695669054
$ws.Cells.Item(4, 30).Value = 117.2
$ws.Cells.Item(4, 31).Value = 4.933187048892024
$ws.Cells.Item(4, 32).Value = 122.133187048892
$ws.Cells.Item(4, 33).Value = 82.13318704889203
$ws.Cells.Item(4, 34).Value = 0.3604522570903564
$ws.Cells.Item(4, 35).Value = 0.4713143400881668
$ws.Cells.Item(4, 36).Value = 0.2748462707907146
$ws.Cells.Item(4, 37).Value = 0.3748094396608526
$ws.Cells.Item(4, 40).Value = 51.40350877192983
$ws.Cells.Item(4, 42).Value = 36.02332765302282

# Row 5
$ws.Cells.Item(5, 9).Value = 0.004657106388284061
$ws.Cells.Item(5, 10).Value = 0.004657106388284061
$ws.Cells.Item(5, 11).Value = 88.7
$ws.Cells.Item(5, 12).Value = 1.190604026845638
$ws.Cells.Item(5, 13).Value = 107.5
$ws.Cells.Item(5, 14).Value = 0.04559141609058908
$ws.Cells.Item(5, 15).Value = 1.211950394588501
$ws.Cells.Item(5, 16).Value = 107.5
$ws.Cells.Item(5, 17).Value = 0.04559141609058908
$ws.Cells.Item(5, 18).Value = 1.211950394588501
$ws.Cells.Item(5, 21).Value = 26.1
$ws.Cells.Item(5, 22).Value = 0.01106917172059884
$ws.Cells.Item(5, 23).Value = 0.06228932584269663
$ws.Cells.Item(5, 24).Value = 0.073535515191195
$ws.Cells.Item(5, 25).Value = -0.01124618934849837
$ws.Cells.Item(5, 26).Value = 0.01921164571924158
$ws.Cells.Item(5, 27).Value = 0.000089470678008530106517
$ws.Cells.Item(5, 28).Value = 0.05192267154871194
$ws.Cells.Item(5, 29).Value = -0.05183320087070341
$ws.Cells.Item(5, 30).Value = 2547.7
$ws.Cells.Item(5, 31).Value = 1.165227870364187
$ws.Cells.Item(5, 32).Value = 2548.865227870364
$ws.Cells.Item(5, 33).Value = 2522.765227870364
$ws.Cells.Item(5, 34).Value = 0.5194593809772764
$ws.Cells.Item(5, 35).Value = 0.6853510416459978
$ws.Cells.Item(5, 36).Value = 0.5168896267386793
$ws.Cells.Item(5, 37).Value = 0.6831272628378297
$ws.Cells.Item(5, 40).Value = 4392.586206896552
$ws.Cells.Item(5, 42).Value = 4349.595220466145
